$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.562.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.193.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.20'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.03%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.586'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.13'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0918'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.529.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.27'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.181.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.775'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.464.14'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.57%  '
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('E23').Value = '  +7.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.37'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.27%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.37%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +13.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0795'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.16'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.107'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.36'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.38%  '
$ws.Range('E38').Value = '  +9.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.15'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.07'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.197'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('B44').Value = 'WOONetwork'
$ws.Range('C44').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.481'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +24.96%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0977'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.09'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('E50').Value = '  +1.88%  '
